# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from 2023-10-08 (45207) to 2023-10-09 (45208).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
